# Fifteenth commit with Customer Page test started
# 1) Dashboard_Tests: remove the "Verify Logout Functionality" test case (rows 19-20)
$wb = $excel.ActiveWorkbook
$dash = $wb.Worksheets.Item("Dashboard_Tests")
$dash.Range("A19:J20").ClearContents()

# Dashboard_Tests is no longer the active tab; reset its selection to C1
$dash.Range("C1").Select()

# 2) Add the new "Customer_Tests" sheet after Dashboard_Tests (last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cust = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$cust.Name = "Customer_Tests"

# Column widths (approximate character widths matching source formatting)
$cust.Columns.Item(1).ColumnWidth = 24.666666666666668
$cust.Columns.Item(2).ColumnWidth = 52.333333333333336
$cust.Columns.Item(3).ColumnWidth = 96.66666666666667

# Header row - copy the header formatting (fill/style) from Auth_Tests!A1:C1
$auth = $wb.Worksheets.Item("Auth_Tests")
$auth.Range("A1:C1").Copy()
$cust.Range("A1:C1").PasteSpecial(-4122)

$cust.Range("A1").Value = "Test Case ID(s)"
$cust.Range("B1").Value = "Test Case Description"
$cust.Range("C1").Value = "Test Steps"

# Test case TC_CUST_01 - Navigate to Customers
$cust.Range("A2").Value = "TC_CUST_01"
$cust.Range("B2").Value = "Navigate to Customers"
$cust.Range("C2").Value = "1.Click on ""Customers"" at ""//a[contains(@href, '/customers')]"""
$cust.Range("C3").Value = "2.Verify URL contains ""customers"""

# Test case TC_CUST_02 - Verify List Consistency
$cust.Range("A4").Value = "TC_CUST_02"
$cust.Range("B4").Value = "Verify List Consistency"
$cust.Range("C4").Value = "1.Verify text ""{DB_QUERY}SELECT COUNT(*) FROM customers"" at ""//table/tbody/tr"""

$cust.Range("C4").Select()
$cust.Activate()
